$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume table is stored as plain text (e.g. "30.125.23",
# "0.07800") where trailing zeros and the exotic "thousands-dot" price
# notation are meaningful. Force the D (Price) and E (Volume) columns to
# a text number format first so Excel does not reinterpret the updated
# values as numbers/dates and strip significant trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.125.23"
$ws.Range("D3").Value = "1.860.12"
$ws.Range("E3").Value = "  -4.60%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "233.47"
$ws.Range("E5").Value = "  -3.98%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "0.4659"
$ws.Range("E7").Value = "  -3.50%  "
$ws.Range("D8").Value = "0.2808"
$ws.Range("E8").Value = "  -3.98%  "
$ws.Range("E9").Value = "  -4.28%  "
$ws.Range("D10").Value = "19.61"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").Value = "0.07800"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "96.44"
$ws.Range("E12").Value = "  -8.50%  "
$ws.Range("D13").Value = "1.863.03"
$ws.Range("E13").Value = "  -4.52%  "
$ws.Range("D14").Value = "5.131"
$ws.Range("E14").Value = "  -3.96%  "
$ws.Range("D15").Value = "0.6654"
$ws.Range("E15").Value = "  -4.21%  "
$ws.Range("D16").Value = "280.57"
$ws.Range("E16").Value = "  -6.10%  "
$ws.Range("D17").Value = "30.155.92"
$ws.Range("E17").Value = "  -3.84%  "
$ws.Range("D18").Value = "0.9998"
$ws.Range("D19").Value = "5.507"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").Value = "12.58"
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").Value = "2.101.01"
$ws.Range("E21").Value = "  -5.11%  "
$ws.Range("D22").Value = "0.000007224"
$ws.Range("E22").Value = "  -5.46%  "
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "6.118"
$ws.Range("E24").Value = "  -5.78%  "
$ws.Range("D25").Value = "9.316"
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("D26").Value = "165.60"
$ws.Range("E26").Value = "  -2.11%  "
$ws.Range("D27").Value = "18.84"
$ws.Range("E27").Value = "  -5.72%  "
$ws.Range("D28").Value = "1.910"
$ws.Range("E28").Value = "  -11.17%  "
$ws.Range("D29").Value = "1.337"
$ws.Range("E29").Value = "  -4.05%  "
$ws.Range("D30").Value = "0.09534"
$ws.Range("E30").Value = "  -6.54%  "
$ws.Range("D31").Value = "4.399"
$ws.Range("E31").Value = "  -5.67%  "
$ws.Range("D32").Value = "1.467"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("D33").Value = "4.097"
$ws.Range("E33").Value = "  -6.52%  "
$ws.Range("D34").Value = "0.04650"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").Value = "0.7007"
$ws.Range("E35").Value = "  -6.37%  "
$ws.Range("D36").Value = "1.090"
$ws.Range("E36").Value = "  -4.37%  "
$ws.Range("D37").Value = "2.703"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("D38").Value = "0.01854"
$ws.Range("E38").Value = "  -5.93%  "
$ws.Range("D39").Value = "6.290"
$ws.Range("E39").Value = "  -4.98%  "
$ws.Range("D40").Value = "2.511"
$ws.Range("E40").Value = "  -5.33%  "
$ws.Range("D41").Value = "73.10"
$ws.Range("E41").Value = "  -5.53%  "
$ws.Range("D42").Value = "0.8517"
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("D43").Value = "1.920"
$ws.Range("E43").Value = "  -6.47%  "
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "103.65"
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("D46").Value = "0.4148"
$ws.Range("E46").Value = "  -5.86%  "
$ws.Range("D47").Value = "993.87"
$ws.Range("E47").Value = "  -3.39%  "
$ws.Range("D48").Value = "7.176"
$ws.Range("E48").Value = "  -6.03%  "
$ws.Range("D49").Value = "9.271"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").Value = "34.08"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").Value = "0.1138"
$ws.Range("E51").Value = "  -6.74%  "
